$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.910.70"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.52"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.10"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.90"
$ws.Range("E8").Value = "  +5.95%  "
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0982"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.110.04"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.36"
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.842.23"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.923.19"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.86"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.06"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.86"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.40"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.123"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.62"
$ws.Range("E29").Value = "  +7.98%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0550"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("E35").Value = "  +11.08%  "
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.742"
$ws.Range("E37").Value = "  +8.50%  "
$ws.Range("E38").Value = "  +11.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "89.78"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0197"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.340.30"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.54"
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.31"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.31"
$ws.Range("E48").Value = "  +75.26%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.025.75"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.39"
$ws.Range("E51").Value = "  +16.51%  "
